# Additional companies sent for questionaire
# Remove the "Location County/City" (col E) and "Parent company" (col B)
# columns from the locomotive list sheet. Deleting from right-to-left so
# column letters stay valid as each delete shifts everything left.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E:E").EntireColumn.Delete()
$ws.Range("B:B").EntireColumn.Delete()

[void]$ws.Range("N9").Select()
